# Generate Report for Handoff
# Updates the localization-status report: flips the "In Translation"
# status to "Ready for handoff" and refreshes the handoff timestamps
# to reflect the newly generated report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Status column(s): "In Translation" -> "Ready for handoff"
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Latest Handoff Datetime / Latest Handoff Date refresh
$overview.Range("D2").Value = "2016-03-23 18:39:58"
$dede.Range("E2").Value = "2016-03-23 18:39:58"

$zhcn.Range("E2").Value = "2016-03-23 18:39:54"
